# NYPD 112th Precinct CompStat weekly report - "New crime data collected"
# Updates the report heading (volume/week-of date) and the Crime Complaints
# table (rows 15-33) with the newly collected weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmtCount = "#,##0"
$fmtPct   = '#,##0.0;"-"#,##0.0'
$fmtGen   = "General"

# ---------------------------------------------------------------------
# Header: issue number and report week
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/27/2024  Through  6/2/2024"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("C15").NumberFormat = $fmtCount
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 150
$ws.Range("L15").Value = -16.666666666666
$ws.Range("M15").Value = 400
$ws.Range("N15").Value = 150

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = -3.030303030303
$ws.Range("L16").Value = -3.030303030303
$ws.Range("M16").Value = -25.581395348837
$ws.Range("N16").Value = -89.003436426116

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").NumberFormat = $fmtGen
$ws.Range("C17").Value = "'0"
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -100
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -54.545454545454
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 44
$ws.Range("K17").Value = -9.090909090909
$ws.Range("L17").Value = -2.439024390243
$ws.Range("M17").Value = 81.818181818181
$ws.Range("N17").Value = -24.528301886792

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 39
$ws.Range("J18").Value = 55
$ws.Range("K18").Value = -29.090909090909
$ws.Range("L18").Value = -13.333333333333
$ws.Range("M18").Value = -20.408163265306
$ws.Range("N18").Value = -93.035714285714

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = 7.894736842105
$ws.Range("I19").Value = 173
$ws.Range("J19").Value = 201
$ws.Range("K19").Value = -13.930348258706
$ws.Range("L19").Value = -23.111111111111
$ws.Range("M19").Value = 16.107382550335
$ws.Range("N19").Value = -57.389162561576

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 13
$ws.Range("E20").Value = 225
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 64.285714285714
$ws.Range("I20").Value = 70
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 40
$ws.Range("L20").Value = 105.882352941176
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = -95.172413793103

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 2.564102564102
$ws.Range("I21").Value = 360
$ws.Range("J21").Value = 385
$ws.Range("K21").Value = -6.493506493506
$ws.Range("L21").Value = -6.735751295336
$ws.Range("M21").Value = 14.649681528662
$ws.Range("N21").Value = -86.975397973950

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("C22").NumberFormat = $fmtGen
$ws.Range("C22").Value = "'0"
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 300

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -30.303030303030
$ws.Range("F24").Value = 111
$ws.Range("G24").Value = 113
$ws.Range("H24").Value = -1.769911504424
$ws.Range("I24").Value = 664
$ws.Range("J24").Value = 674
$ws.Range("K24").Value = -1.483679525222
$ws.Range("L24").Value = -13.089005235602
$ws.Range("M24").Value = 66.834170854271

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 28
$ws.Range("E25").Value = -35.714285714285
$ws.Range("F25").Value = 76
$ws.Range("G25").Value = 79
$ws.Range("H25").Value = -3.797468354430
$ws.Range("I25").Value = 495
$ws.Range("J25").Value = 480
$ws.Range("K25").Value = 3.125
$ws.Range("L25").Value = -12.078152753108

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 12.5
$ws.Range("I26").Value = 120
$ws.Range("J26").Value = 104
$ws.Range("K26").Value = 15.384615384615
$ws.Range("L26").Value = 53.846153846153
$ws.Range("M26").Value = 29.032258064516

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("C27").NumberFormat = $fmtCount
$ws.Range("C27").Value = 1
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 60
$ws.Range("L27").Value = -11.111111111111

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C28").Value = 3
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 13
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -18.75

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------
$ws.Range("D31").Value = 1
$ws.Range("G31").Value = 3
$ws.Range("J31").Value = 12
$ws.Range("K31").Value = -58.333333333333

# ---------------------------------------------------------------------
# Row 33 - Traffic Fatalities
# ---------------------------------------------------------------------
$ws.Range("D33").NumberFormat = $fmtCount
$ws.Range("D33").Value = 1
$ws.Range("E33").NumberFormat = $fmtPct
$ws.Range("E33").Value = -100
$ws.Range("G33").NumberFormat = $fmtCount
$ws.Range("G33").Value = 1
$ws.Range("H33").NumberFormat = $fmtPct
$ws.Range("H33").Value = 0
$ws.Range("J33").NumberFormat = $fmtCount
$ws.Range("J33").Value = 1
$ws.Range("K33").NumberFormat = $fmtPct
$ws.Range("K33").Value = 0
